$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 already exists (empty, styled with quote-prefix text style) - set its value
# Leading apostrophe forces Excel to store the value as text (quote-prefixed)
$ws.Range("A3").Value = "'8939465567"

# A4 is a new cell - set its value
$ws.Range("A4").Value = "'8979466578"

# Update the selection to match the new active cell
$ws.Range("A4").Select()
